$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 67  # E2
$ws.Cells.Item(2, 6).Value = 38  # F2
$ws.Cells.Item(2, 8).Value = 38  # H2

$ws.Cells.Item(5, 5).Value = 73  # E5
$ws.Cells.Item(5, 6).Value = 39  # F5
$ws.Cells.Item(5, 8).Value = 39  # H5

$ws.Cells.Item(10, 5).Value = 218  # E10
$ws.Cells.Item(10, 6).Value = 94  # F10
$ws.Cells.Item(10, 8).Value = 94  # H10

$ws.Cells.Item(11, 5).Value = 167  # E11
$ws.Cells.Item(11, 6).Value = 81  # F11
$ws.Cells.Item(11, 8).Value = 81  # H11

$ws.Cells.Item(12, 5).Value = 244  # E12
$ws.Cells.Item(12, 6).Value = 121  # F12
$ws.Cells.Item(12, 8).Value = 121  # H12

$ws.Cells.Item(13, 5).Value = 78  # E13

$ws.Cells.Item(14, 5).Value = 71  # E14
$ws.Cells.Item(14, 6).Value = 29  # F14
$ws.Cells.Item(14, 8).Value = 29  # H14

$ws.Cells.Item(15, 5).Value = 98  # E15

$ws.Cells.Item(16, 5).Value = 96  # E16
$ws.Cells.Item(16, 6).Value = 45  # F16
$ws.Cells.Item(16, 8).Value = 45  # H16

$ws.Cells.Item(17, 6).Value = 20  # F17
$ws.Cells.Item(17, 8).Value = 20  # H17

$ws.Cells.Item(18, 6).Value = 14  # F18
$ws.Cells.Item(18, 8).Value = 14  # H18

$ws.Cells.Item(20, 5).Value = 59  # E20
$ws.Cells.Item(20, 6).Value = 21  # F20
$ws.Cells.Item(20, 8).Value = 21  # H20

$ws.Cells.Item(21, 5).Value = 78  # E21

$ws.Cells.Item(22, 5).Value = 95  # E22

$ws.Cells.Item(23, 5).Value = 105  # E23

$ws.Cells.Item(24, 5).Value = 111  # E24
$ws.Cells.Item(24, 6).Value = 52  # F24
$ws.Cells.Item(24, 8).Value = 52  # H24

$ws.Cells.Item(25, 5).Value = 109  # E25
$ws.Cells.Item(25, 6).Value = 41  # F25
$ws.Cells.Item(25, 8).Value = 41  # H25

$ws.Cells.Item(26, 5).Value = 66  # E26
$ws.Cells.Item(26, 6).Value = 30  # F26
$ws.Cells.Item(26, 8).Value = 30  # H26

$ws.Cells.Item(27, 5).Value = 161  # E27
$ws.Cells.Item(27, 6).Value = 77  # F27
$ws.Cells.Item(27, 8).Value = 77  # H27

$ws.Cells.Item(28, 5).Value = 96  # E28
$ws.Cells.Item(28, 6).Value = 27  # F28
$ws.Cells.Item(28, 8).Value = 27  # H28

$ws.Cells.Item(29, 5).Value = 105  # E29

$ws.Cells.Item(30, 5).Value = 116  # E30
$ws.Cells.Item(30, 6).Value = 56  # F30
$ws.Cells.Item(30, 8).Value = 56  # H30

$ws.Cells.Item(32, 5).Value = 110  # E32
$ws.Cells.Item(32, 6).Value = 51  # F32
$ws.Cells.Item(32, 8).Value = 51  # H32

$ws.Cells.Item(33, 5).Value = 148  # E33
$ws.Cells.Item(33, 6).Value = 66  # F33
$ws.Cells.Item(33, 8).Value = 66  # H33

$ws.Cells.Item(34, 5).Value = 115  # E34
$ws.Cells.Item(34, 6).Value = 57  # F34
$ws.Cells.Item(34, 8).Value = 57  # H34

$ws.Cells.Item(35, 5).Value = 73  # E35
$ws.Cells.Item(35, 6).Value = 36  # F35
$ws.Cells.Item(35, 8).Value = 36  # H35

$ws.Cells.Item(38, 5).Value = 52  # E38

$ws.Cells.Item(39, 5).Value = 116  # E39

$ws.Cells.Item(40, 5).Value = 146  # E40
$ws.Cells.Item(40, 6).Value = 58  # F40
$ws.Cells.Item(40, 8).Value = 58  # H40

$ws.Cells.Item(41, 5).Value = 196  # E41
$ws.Cells.Item(41, 6).Value = 69  # F41
$ws.Cells.Item(41, 8).Value = 69  # H41

$ws.Cells.Item(42, 5).Value = 176  # E42
$ws.Cells.Item(42, 6).Value = 83  # F42
$ws.Cells.Item(42, 8).Value = 83  # H42

$ws.Cells.Item(43, 5).Value = 57  # E43
$ws.Cells.Item(43, 6).Value = 23  # F43
$ws.Cells.Item(43, 8).Value = 23  # H43

$ws.Cells.Item(44, 5).Value = 152  # E44
$ws.Cells.Item(44, 6).Value = 68  # F44
$ws.Cells.Item(44, 8).Value = 68  # H44

$ws.Cells.Item(45, 5).Value = 65  # E45

$ws.Cells.Item(46, 5).Value = 139  # E46
$ws.Cells.Item(46, 6).Value = 59  # F46
$ws.Cells.Item(46, 8).Value = 59  # H46

$ws.Cells.Item(47, 5).Value = 228  # E47
$ws.Cells.Item(47, 6).Value = 90  # F47
$ws.Cells.Item(47, 8).Value = 90  # H47

$ws.Cells.Item(48, 5).Value = 113  # E48

$ws.Cells.Item(49, 5).Value = 124  # E49
$ws.Cells.Item(49, 6).Value = 52  # F49
$ws.Cells.Item(49, 8).Value = 52  # H49

$ws.Cells.Item(50, 5).Value = 104  # E50
$ws.Cells.Item(50, 6).Value = 37  # F50
$ws.Cells.Item(50, 8).Value = 37  # H50

$ws.Cells.Item(52, 5).Value = 8  # E52
$ws.Cells.Item(52, 6).Value = 5  # F52
$ws.Cells.Item(52, 8).Value = 5  # H52

